$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New region: "Чойский муниципальный район" (Altai) - new local extremum
# Two rows: female (row 42) and male (row 43)

$ws.Range("A42:R43").HorizontalAlignment = -4108

$ws.Range("A42").Value = 84645000
$ws.Range("B42").Value = "Чойский муниципальный район"
$ws.Range("C42").Value = "female"
$ws.Range("D42").Value = 2018
$ws.Range("E42").Value = 0.0715
$ws.Range("F42").Value = 0.09247
$ws.Range("G42").Value = 0.0768
$ws.Range("H42").Value = 0.04877
$ws.Range("I42").Value = 0.0396
$ws.Range("J42").Value = 0.035
$ws.Range("K42").Value = 0.0855
$ws.Range("L42").Value = 0.06885
$ws.Range("M42").Value = 0.0739
$ws.Range("N42").Value = 0.0654
$ws.Range("O42").Value = 0.1405
$ws.Range("P42").Value = 0.07947
$ws.Range("Q42").Value = 0.0751
$ws.Range("R42").Value = 0.0471

$ws.Range("A43").Value = 84645000
$ws.Range("B43").Value = "Чойский муниципальный район"
$ws.Range("C43").Value = "male"
$ws.Range("D43").Value = 2018
$ws.Range("E43").Value = 0.0775
$ws.Range("F43").Value = 0.0945
$ws.Range("G43").Value = 0.07635
$ws.Range("H43").Value = 0.06216
$ws.Range("I43").Value = 0.0475
$ws.Range("J43").Value = 0.04028
$ws.Range("K43").Value = 0.0752
$ws.Range("L43").Value = 0.0789
$ws.Range("M43").Value = 0.07263
$ws.Range("N43").Value = 0.0554
$ws.Range("O43").Value = 0.1276
$ws.Range("P43").Value = 0.0787
$ws.Range("Q43").Value = 0.0715
$ws.Range("R43").Value = 0.0419

$ws.Range("B45").Select() | Out-Null
